$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) Insert 4 blank rows above the old summary block (old row 134 -> row 138).
#    Formulas referencing the summary cells (C134/E134/G134/...) shift
#    automatically to the new row numbers when Excel inserts whole rows.
# ---------------------------------------------------------------------------
$ws.Range("A134:A137").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2) New "start time" formulas that chain to the previous session's end time.
# ---------------------------------------------------------------------------
$ws.Range("J127").Formula = "=K126"
$ws.Range("J129").Formula = "=K128"

# ---------------------------------------------------------------------------
# 3) Turn the previously-blank row 131 into a new tracked task entry
#    (continuation of Issue 22 / "Interface Design" / "MockUps"),
#    for the new friend-request screen feature.
# ---------------------------------------------------------------------------
$ws.Range("A131").Value = 22
$ws.Range("B131").Value = "Interface Design"
$ws.Range("C131").Value = "MockUps"
$ws.Range("D131").Value = "[FEATURE]"
$ws.Range("E131").Value = "Freundschaftsanfrage ausstehend"
$ws.Range("F131").Value = 44460
$ws.Range("G131").Value = 44481
$ws.Range("J131").Formula = "=K130"
$ws.Range("K131").Value = 0.47916666666666669
$ws.Range("I131").Formula = "=ROUNDUP(((SUM(K131-J131)*24*60/60)/0.25),0)*0.25"

# Apply the same visual styling used by the rows above (copy number formats /
# fonts / alignment only, so the values/formulas just written are preserved).
$ws.Range("A126:G126").Copy()
$ws.Range("A131:G131").PasteSpecial(-4122)
$ws.Range("I126:K126").Copy()
$ws.Range("I131:K131").PasteSpecial(-4122)

# Re-apply values/formulas touched indirectly by PasteSpecial (format paste
# does not alter cell contents, but keep these explicit/authoritative).
$ws.Range("D131").Value = "[FEATURE]"
$ws.Range("J131").Formula = "=K130"
$ws.Range("K131").Value = 0.47916666666666669
$ws.Range("I131").Formula = "=ROUNDUP(((SUM(K131-J131)*24*60/60)/0.25),0)*0.25"

# ---------------------------------------------------------------------------
# 4) Rows 132/133 gain the same (empty) styled A/B/C/I/J cells as row 131.
# ---------------------------------------------------------------------------
$ws.Range("A126:C126").Copy()
$ws.Range("A132:C133").PasteSpecial(-4122)
$ws.Range("I126").Copy()
$ws.Range("I132:I133").PasteSpecial(-4122)
$ws.Range("J126").Copy()
$ws.Range("J132:J133").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) Extend the prefix data-validation range to cover the newly inserted rows
#    (D115:D133 -> D115:D137). Only the touched area is removed/recreated so
#    the untouched D2:D3/D13:D17/D22:D27/D35:D36/D41:D113 portion keeps its
#    original single dataValidation entry.
# ---------------------------------------------------------------------------
$ws.Range("D115:D133").Validation.Delete()
$dv = $ws.Range("D115:D137")
$dv.Validation.Add(3, 1, 1, '=$N$3:$N$6')
$dv.Validation.IgnoreBlank = $true
$dv.Validation.InCellDropdown = $true
$dv.Validation.ShowInput = $true
$dv.Validation.ShowError = $true
$dv.Validation.ErrorTitle = "Prefix nicht unterstützt"
$dv.Validation.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden_x000a_"
$dv.Validation.InputTitle = "Prefix"
$dv.Validation.InputMessage = "Wählen Sie einen Prefix aus"

# ---------------------------------------------------------------------------
# 6) Restore the sheet selection / scroll position to match the new layout.
# ---------------------------------------------------------------------------
$ws.Range("A120:I141").Select()
$ws.Range("I141").Activate()

Write-Host "done"
